$p = $ppt.ActivePresentation

# --- Slide 15: "Available Datasets & Source Code" ----------------------
# Update the trailing "GitHub: " paragraph into a bold "Github Repo: "
# label plus a live hyperlink to the project's GitHub repository.
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(7)

$startPos = $para.Start
$url = "https://github.com/kotlarmilos/anomalydetection"
$fullText = "Github Repo: " + $url

# Rewrite the paragraph's text (diff-preserves the run's lang/dirty rPr).
$para.Text = $fullText

# "Github" -> bold
$runGithub = $tr.Characters($startPos, 6)
$runGithub.Font.Bold = 1

# " Repo" -> bold
$runRepo = $tr.Characters($startPos + 6, 5)
$runRepo.Font.Bold = 1

# ": " stays regular (no formatting change needed)

# URL run -> hyperlink to the GitHub repo
$urlStart = $startPos + 13
$runUrl = $tr.Characters($urlStart, $url.Length)
$runUrl.ActionSettings(1).Hyperlink.Address = $url
